$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5262118
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 5815972.5
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 17447917.5
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -17448253.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3584.05
$ws.Range("I32").Value = 3423.0737
$ws.Range("J32").Value = 6642.6
$ws.Range("K32").Value = 3423.0737
$ws.Range("L32").Value = 6642.6
$ws.Range("M32").Value = -3136.0737
$ws.Range("N32").Value = -7216.6

$ws.Range("H61").Value = 2049.3403
$ws.Range("I61").Value = 2093.2927
$ws.Range("K61").Value = 2093.2927
$ws.Range("M61").Value = -1881.2927

$ws.Range("H136").Value = 2049.3403
$ws.Range("I136").Value = 2093.2927
$ws.Range("K136").Value = 6279.8781
$ws.Range("M136").Value = -3729.8781

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3977.8125
$ws.Range("I105").Value = 1603.3334
$ws.Range("J105").Value = 4223.448
$ws.Range("K105").Value = 1603.3334
$ws.Range("L105").Value = 4223.448
$ws.Range("M105").Value = 143.6666
$ws.Range("N105").Value = -7717.448

$ws.Range("H134").Value = 3679.6785
$ws.Range("I134").Value = 2537.7576
$ws.Range("J134").Value = 5318.087
$ws.Range("K134").Value = 7613.2728
$ws.Range("L134").Value = 15954.261
$ws.Range("M134").Value = -5078.2728
$ws.Range("N134").Value = -21024.261

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3116
$ws.Range("I16").Value = 2466
$ws.Range("J16").Value = 3766
$ws.Range("K16").Value = 2466
$ws.Range("L16").Value = 3766
$ws.Range("M16").Value = -2179
$ws.Range("N16").Value = -4340

$ws.Range("H31").Value = 2930.0378
$ws.Range("I31").Value = 2626.1155
$ws.Range("J31").Value = 3222.7036
$ws.Range("K31").Value = 2626.1155
$ws.Range("L31").Value = 3222.7036
$ws.Range("M31").Value = -2331.1155
$ws.Range("N31").Value = -3812.7036

$ws.Range("H34").Value = 2930.0378
$ws.Range("I34").Value = 2626.1155
$ws.Range("J34").Value = 3222.7036
$ws.Range("K34").Value = 2626.1155
$ws.Range("L34").Value = 3222.7036
$ws.Range("M34").Value = -2424.1155
$ws.Range("N34").Value = -3626.7036

$ws.Range("H113").Value = 3116
$ws.Range("I113").Value = 2466
$ws.Range("J113").Value = 3766
$ws.Range("K113").Value = 2466
$ws.Range("L113").Value = 3766
$ws.Range("M113").Value = -296
$ws.Range("N113").Value = -8106

$ws.Range("H122").Value = 1797.8
$ws.Range("I122").Value = 1884.8889
$ws.Range("J122").Value = 1014
$ws.Range("K122").Value = 5654.6667
$ws.Range("L122").Value = 3042
$ws.Range("M122").Value = -3204.6667
$ws.Range("N122").Value = -7942

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 34483500
$ws.Range("I113").Value = 142857700
$ws.Range("J113").Value = 800.63635
$ws.Range("K113").Value = 428573100
$ws.Range("L113").Value = 2401.90905
$ws.Range("M113").Value = -428570930
$ws.Range("N113").Value = -6741.90905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3032.6
$ws.Range("I102").Value = 3830.7368
$ws.Range("J102").Value = 1654
$ws.Range("K102").Value = 3830.7368
$ws.Range("L102").Value = 1654
$ws.Range("M102").Value = -2208.7368
$ws.Range("N102").Value = -4898

$ws.Range("H122").Value = 2208.1667
$ws.Range("I122").Value = 2057.5386
$ws.Range("J122").Value = 2599.8
$ws.Range("K122").Value = 6172.6158
$ws.Range("L122").Value = 7799.400000000001
$ws.Range("M122").Value = -3722.6158
$ws.Range("N122").Value = -12699.4

$ws.Range("H132").Value = 5885.125
$ws.Range("I132").Value = 6712.4165
$ws.Range("J132").Value = 3403.25
$ws.Range("K132").Value = 20137.2495
$ws.Range("L132").Value = 10209.75
$ws.Range("M132").Value = -17607.2495
$ws.Range("N132").Value = -15269.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 654.7308
$ws.Range("I46").Value = 551.8
$ws.Range("J46").Value = 795.0909
$ws.Range("K46").Value = 551.8
$ws.Range("L46").Value = 795.0909
$ws.Range("M46").Value = -363.8
$ws.Range("N46").Value = -1171.0909

$ws.Range("H55").Value = 240.47826
$ws.Range("I55").Value = 105
$ws.Range("J55").Value = 299.75
$ws.Range("K55").Value = 105
$ws.Range("L55").Value = 299.75
$ws.Range("M55").Value = 68
$ws.Range("N55").Value = -645.75

$ws.Range("H61").Value = 1940.6154
$ws.Range("I61").Value = 1208.3334
$ws.Range("J61").Value = 2568.2856
$ws.Range("K61").Value = 1208.3334
$ws.Range("L61").Value = 2568.2856
$ws.Range("M61").Value = -1006.3334
$ws.Range("N61").Value = -2972.2856

$ws.Range("H81").Value = 37593.668
$ws.Range("J81").Value = 37593.668
$ws.Range("L81").Value = 37593.668
$ws.Range("N81").Value = -39589.668

$ws.Range("H84").Value = 37593.668
$ws.Range("J84").Value = 37593.668
$ws.Range("L84").Value = 112781.004
$ws.Range("N84").Value = -122765.004

$ws.Range("H113").Value = 1940.6154
$ws.Range("I113").Value = 1208.3334
$ws.Range("J113").Value = 2568.2856
$ws.Range("K113").Value = 1208.3334
$ws.Range("L113").Value = 2568.2856
$ws.Range("M113").Value = 961.6666
$ws.Range("N113").Value = -6908.2856

$ws.Range("H122").Value = 7132.696
$ws.Range("I122").Value = 6057.846
$ws.Range("J122").Value = 8530
$ws.Range("K122").Value = 18173.538
$ws.Range("L122").Value = 25590
$ws.Range("M122").Value = -15723.538
$ws.Range("N122").Value = -30490

$ws.Range("H132").Value = 7698831
$ws.Range("I132").Value = 4737.549
$ws.Range("J132").Value = 35727316
$ws.Range("K132").Value = 14212.647
$ws.Range("L132").Value = 107181948
$ws.Range("M132").Value = -11682.647
$ws.Range("N132").Value = -107187008

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 5250
$ws.Range("I113").Value = 10000
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 30000
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -27830
$ws.Range("N113").Value = -5840

$ws.Range("H122").Value = 2401.6765
$ws.Range("I122").Value = 2332.8965
$ws.Range("J122").Value = 2800.6
$ws.Range("K122").Value = 6998.689499999999
$ws.Range("L122").Value = 8401.799999999999
$ws.Range("M122").Value = -4548.689499999999
$ws.Range("N122").Value = -13301.8

$ws.Range("H132").Value = 1360.228
$ws.Range("I132").Value = 1179.3077
$ws.Range("J132").Value = 3241.8
$ws.Range("K132").Value = 3537.9231
$ws.Range("L132").Value = 9725.400000000001
$ws.Range("M132").Value = -1007.9231
$ws.Range("N132").Value = -14785.4
